$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a text value into a cell without Excel's "smart" value
# parser reinterpreting look-alike text (e.g. "$000") as a currency
# number. Forcing a text NumberFormat before the write, then resetting
# the cell back to the default "Normal" style after, gives us a plain
# shared-string cell with no stray style index - matching how the
# source file stores these note labels.
# ---------------------------------------------------------------------
function Set-TextCell($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# =====================================================================
# Section 1: restore the missing "Total" keyword on note rows whose
# label lost its trailing value (commit: "added total keywords to note
# df where particular value is not present").
# =====================================================================
$totalEdits = @(
    @{ Sheet = 1;  Cells = @{ "A4" = " Total"; "A7" = " Total" } },
    @{ Sheet = 2;  Cells = @{
            "A2"  = " Total"
            "A5"  = "Related party receivables: Total"
            "A8"  = "Related party receivables: Total"
            "A9"  = "Related party receivables: Total"
            "A12" = "Non-current Total"
            "A13" = " Total"
            "A16" = "Related party receivables: Total"
            "A19" = "Related party receivables: Total"
            "A20" = "Related party receivables: Total"
            "A23" = "Non-current Total"
        }
    },
    @{ Sheet = 5;  Cells = @{ "A4" = "Current Total"; "A7" = "Current Total" } },
    @{ Sheet = 11; Cells = @{
            "A4"  = " Total"
            "A7"  = "Related party payables: Total"
            "A8"  = "Related party payables: Total"
            "A11" = " Total"
            "A14" = "Related party payables: Total"
            "A15" = "Related party payables: Total"
        }
    },
    @{ Sheet = 14; Cells = @{
            "A5"  = "Current Total"
            "A9"  = "Non-current Total"
            "A14" = "Current Total"
            "A18" = "Non-current Total"
        }
    },
    @{ Sheet = 16; Cells = @{ "A2" = "Total 12,540,000 ordinary shares fully issued and paid 2020.12,540,000)" } },
    @{ Sheet = 19; Cells = @{ "A6" = " Total"; "A11" = " Total" } },
    @{ Sheet = 20; Cells = @{ "A4" = " Total"; "A7" = " Total" } },
    @{ Sheet = 21; Cells = @{ "A5" = " Total"; "A9" = " Total" } },
    @{ Sheet = 23; Cells = @{ "A4" = " Total"; "A7" = " Total" } }
)

foreach ($edit in $totalEdits) {
    $ws = $wb.Worksheets.Item($edit.Sheet)
    foreach ($ref in $edit.Cells.Keys) {
        Set-TextCell $ws.Range($ref) $edit.Cells[$ref]
    }
}

# =====================================================================
# Section 2: the note table that used to live on the last sheet
# ("4_7_412c019a-25e7-31") grows from 4 data rows to 14 data rows
# (it now also reports the tax-reconciliation detail rows).
# =====================================================================
$sheet24 = $wb.Worksheets.Item($wb.Worksheets.Count)

$sheet24Rows = @(
    @(" Accounting profit before income tax", "$000", 2021, 63330),
    @(" At the statutory income tax rate of 30% (2020: 30%)", "$000", 2021, 18999),
    @(" Adjustments in respect of current income tax of previous years", "$000", 2021, 269),
    @(" Expenditure not allowable for income tax purpose", "$000", 2021, -1708),
    @(" Others", "$000", 2021, -203),
    @(" Income tax expense reported in the consolidated statement of profit or loss", "$000", 2021, 17357),
    @(" Effective income tax rate", "$000", 2021, 0),
    @(" Accounting profit before income tax", "$000", 2020, 43450),
    @(" At the statutory income tax rate of 30% (2020: 30%)", "$000", 2020, 13035),
    @(" Adjustments in respect of current income tax of previous years", "$000", 2020, 311),
    @(" Expenditure not allowable for income tax purpose", "$000", 2020, 21),
    @(" Others", "$000", 2020, 174),
    @(" Income tax expense reported in the consolidated statement of profit or loss", "$000", 2020, 13541),
    @(" Effective income tax rate", "$000", 2020, 0)
)

# clear any previous data rows beyond the header before rewriting
$sheet24.Range("A2:D100").ClearContents()

$r = 2
foreach ($row in $sheet24Rows) {
    Set-TextCell $sheet24.Cells.Item($r, 1) $row[0]
    Set-TextCell $sheet24.Cells.Item($r, 2) $row[1]
    $sheet24.Cells.Item($r, 3).Value = $row[2]
    $sheet24.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# =====================================================================
# Section 3: append a brand-new note sheet, "4_7_d493819c-94d0-30",
# holding the 4-row table that the old sheet used to carry (values
# updated for the "Employee entitlements" note).
# =====================================================================
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $sheet24)
$newSheet.Name = "4_7_d493819c-94d0-30"

Set-TextCell $newSheet.Cells.Item(1, 1) "rows"
Set-TextCell $newSheet.Cells.Item(1, 2) "columns"
Set-TextCell $newSheet.Cells.Item(1, 3) "year"
Set-TextCell $newSheet.Cells.Item(1, 4) "value"
$newSheet.Range("A1:D1").Font.Bold = $true

$sheet25Rows = @(
    @(" Employee entitlements", "Consolidated statement of Consolidated statement of financial position profit or loss $000", 2021, 2339),
    @(" Employee entitlements", "Consolidated statement of Consolidated statement of financial position profit or loss $000", 2020, 2054),
    @(" Employee entitlements", "Consolidated statement of Consolidated statement of financial position profit or loss $000", 2021, -285),
    @(" Employee entitlements", "Consolidated statement of Consolidated statement of financial position profit or loss $000", 2020, 167)
)

$r = 2
foreach ($row in $sheet25Rows) {
    Set-TextCell $newSheet.Cells.Item($r, 1) $row[0]
    Set-TextCell $newSheet.Cells.Item($r, 2) $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
